# Scheduled-runner price/profit refresh: updates cached market-price and
# profit columns (H, I, J, K, L, M, N) on several rows across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR sheets. A few rows also lose a now-blank
# trailing cell (ClearContents) where the source data no longer supplies a
# value for that column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3177.1538
$ws.Range("I62").Value = 1549.6666
$ws.Range("K62").Value = 1549.6666
$ws.Range("M62").Value = -925.6666
$ws.Range("H65").Value = 3177.1538
$ws.Range("I65").Value = 1549.6666
$ws.Range("K65").Value = 7748.333000000001
$ws.Range("M65").Value = -4628.333000000001
$ws.Range("H106").Value = 2621.9048
$ws.Range("I106").Value = 2475.7222
$ws.Range("K106").Value = 2475.7222
$ws.Range("M106").Value = -1844.7222
$ws.Range("H138").Value = 2540.8518
$ws.Range("J138").Value = 2789.2666
$ws.Range("L138").Value = 8367.799800000001
$ws.Range("N138").Value = -18647.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 927768
$ws.Range("I2").Value = 1043015.94
$ws.Range("K2").Value = 1043015.94
$ws.Range("M2").Value = -1042902.94
$ws.Range("H27").Value = 42558.332
$ws.Range("J27").Value = 42558.332
$ws.Range("L27").Value = 42558.332
$ws.Range("N27").Value = -42926.332
$ws.Range("H61").Value = 27030810
$ws.Range("I61").Value = 31253768
$ws.Range("J61").Value = 3881.2
$ws.Range("K61").Value = 31253768
$ws.Range("L61").Value = 3881.2
$ws.Range("M61").Value = -31253556
$ws.Range("N61").Value = -4305.2
$ws.Range("H74").Value = 25004338
$ws.Range("I74").Value = 35719390
$ws.Range("K74").Value = 35719390
$ws.Range("M74").Value = -35718516
$ws.Range("H77").Value = 25004338
$ws.Range("I77").Value = 35719390
$ws.Range("K77").Value = 178596950
$ws.Range("M77").Value = -178592582
$ws.Range("H110").Value = 54745.26
$ws.Range("I110").Value = 68289.92999999999
$ws.Range("K110").Value = 68289.92999999999
$ws.Range("M110").Value = -66244.92999999999
$ws.Range("H116").Value = 927768
$ws.Range("I116").Value = 1043015.94
$ws.Range("K116").Value = 1043015.94
$ws.Range("M116").Value = -1040721.94
$ws.Range("H122").Value = 2144.3333
$ws.Range("I122").Value = 2521.8
$ws.Range("J122").Value = 257
$ws.Range("K122").Value = 7565.400000000001
$ws.Range("L122").Value = 771
$ws.Range("M122").Value = -5115.400000000001
$ws.Range("N122").Value = -5671
$ws.Range("H132").Value = 3129319
$ws.Range("I132").Value = 3129319
$ws.Range("K132").Value = 9387957
$ws.Range("M132").Value = -9385427
$ws.Range("H136").Value = 27030810
$ws.Range("I136").Value = 31253768
$ws.Range("J136").Value = 3881.2
$ws.Range("K136").Value = 93761304
$ws.Range("L136").Value = 11643.6
$ws.Range("M136").Value = -93758754
$ws.Range("N136").Value = -16743.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 927768
$ws.Range("I3").Value = 1043015.94
$ws.Range("K3").Value = 1043015.94
$ws.Range("M3").Value = -1042901.94

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 23194.408
$ws.Range("I60").Value = 17569.5
$ws.Range("J60").Value = 24444.389
$ws.Range("K60").Value = 17569.5
$ws.Range("L60").Value = 24444.389
$ws.Range("M60").Value = -17058.5
$ws.Range("N60").Value = -25466.389
$ws.Range("H99").Value = 1678.25
$ws.Range("I99").Value = 1603
$ws.Range("K99").Value = 1603
$ws.Range("M99").Value = -105
$ws.Range("H122").Value = 1511.7606
$ws.Range("I122").Value = 1486.2623
$ws.Range("J122").Value = 1667.3
$ws.Range("K122").Value = 4458.7869
$ws.Range("L122").Value = 5001.9
$ws.Range("M122").Value = -2008.7869
$ws.Range("N122").Value = -9901.9
$ws.Range("H126").Value = 1678.25
$ws.Range("I126").Value = 1603
$ws.Range("K126").Value = 4809
$ws.Range("M126").Value = -2339
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 151782
$ws.Range("I11").Value = 159759.89
$ws.Range("K11").Value = 479279.67
$ws.Range("M11").Value = -479139.67
$ws.Range("H34").Value = 2503.9167
$ws.Range("J34").Value = 3715.1428
$ws.Range("L34").Value = 11145.4284
$ws.Range("N34").Value = -11313.4284
$ws.Range("H39").Value = 1199.6
$ws.Range("J39").Value = 1999
$ws.Range("L39").Value = 5997
$ws.Range("N39").Value = -6585
$ws.Range("H75").Value = 2126.75
$ws.Range("J75").Value = 2502.3333
$ws.Range("L75").Value = 7506.999899999999
$ws.Range("N75").Value = -9502.999899999999
$ws.Range("H78").Value = 2126.75
$ws.Range("J78").Value = 2502.3333
$ws.Range("L78").Value = 22520.9997
$ws.Range("N78").Value = -32504.9997
$ws.Range("H86").Value = 351
$ws.Range("I86").Value = 351
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1053
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 133
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 351
$ws.Range("I89").Value = 351
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3159
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2769
$ws.Range("N89").ClearContents()
$ws.Range("H131").Value = 2239.1875
$ws.Range("I131").Value = 1851.6
$ws.Range("K131").Value = 5554.799999999999
$ws.Range("M131").Value = -514.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 664.8461
$ws.Range("I107").Value = 466.09525
$ws.Range("J107").Value = 1499.6
$ws.Range("K107").Value = 466.09525
$ws.Range("L107").Value = 1499.6
$ws.Range("M107").Value = 1453.90475
$ws.Range("N107").Value = -5339.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4344.4375
$ws.Range("I7").Value = 4344.4375
$ws.Range("K7").Value = 4344.4375
$ws.Range("M7").Value = -4232.4375
$ws.Range("H46").Value = 776.125
$ws.Range("I46").Value = 789.4286
$ws.Range("J46").Value = 765.7778
$ws.Range("K46").Value = 789.4286
$ws.Range("L46").Value = 765.7778
$ws.Range("M46").Value = -601.4286
$ws.Range("N46").Value = -1141.7778
$ws.Range("H122").Value = 3856.1538
$ws.Range("J122").Value = 3799.6
$ws.Range("L122").Value = 11398.8
$ws.Range("N122").Value = -16298.8
$ws.Range("H126").Value = 4344.4375
$ws.Range("I126").Value = 4344.4375
$ws.Range("K126").Value = 13033.3125
$ws.Range("M126").Value = -10563.3125
$ws.Range("H132").Value = 40010364
$ws.Range("I132").Value = 43647572
$ws.Range("K132").Value = 130942716
$ws.Range("M132").Value = -130940186
$ws.Range("H136").Value = 2411.3044
$ws.Range("I136").Value = 1697.4
$ws.Range("J136").Value = 2960.4614
$ws.Range("K136").Value = 5092.200000000001
$ws.Range("L136").Value = 8881.3842
$ws.Range("M136").Value = -2542.200000000001
$ws.Range("N136").Value = -13981.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H132").Value = 11114866
$ws.Range("I132").Value = 11908571
$ws.Range("K132").Value = 35725713
$ws.Range("M132").Value = -35723183
